# Weekly update: insert a new price record as row 161 (pushing the
# existing rows 161-251 down to 162-252), which is how this "logica_diaria"
# sheet appends newly published daily observations at the top of its
# history window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 161; everything below (old 161..251)
# shifts down to 162..252, and the sheet's dimension grows to A1:R252.
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(161, 1).Value  = 9
$ws.Cells.Item(161, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(161, 3).Value  = "Metropolitana"
$ws.Cells.Item(161, 4).Value  = 44529
$ws.Cells.Item(161, 5).Value  = 13
$ws.Cells.Item(161, 6).Value  = 100112044
$ws.Cells.Item(161, 7).Value  = "Perejil"
$ws.Cells.Item(161, 8).Value  = "Sin especificar"
$ws.Cells.Item(161, 9).Value  = "Primera"
$ws.Cells.Item(161, 10).Value = 61
$ws.Cells.Item(161, 11).Value = 12000
$ws.Cells.Item(161, 12).Value = 13000
$ws.Cells.Item(161, 13).Value = 12492
$ws.Cells.Item(161, 14).Value = "$/docena de atados"
$ws.Cells.Item(161, 15).Value = "Región Metropolitana"
$ws.Cells.Item(161, 16).Value = 4164
$ws.Cells.Item(161, 17).Value = 3
$ws.Cells.Item(161, 18).Value = "Hortaliza"
